$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B3: "ssss" -> "ss2"
$ws.Range("B3").Value = "ss2"

# --- Row 6: id=4, "ssss" across B:G ---
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "ssss"
$ws.Range("C6").Value = "ssss"
$ws.Range("D6").Value = "ssss"
$ws.Range("E6").Value = "ssss"
$ws.Range("F6").Value = "ssss"
$ws.Range("G6").Value = "ssss"

# --- Row 7: id=5, "228" (as text) across B:G ---
$ws.Range("A7").Value = 5

# Plain .Value assignment of an all-digit string like "228" gets
# auto-detected as a number (same as typing it into Excel). To land it as
# literal text (matching t="inlineStr" in the target) without touching
# NumberFormat (which would register a new style not present in the
# target), build it as a text formula in a scratch cell, copy it, and
# paste-special just the resulting value into each destination cell.
$ws.Range("I1").Formula = "=""228"""
$ws.Range("I1").Copy()
$ws.Range("B7").PasteSpecial(-4163)
$ws.Range("C7").PasteSpecial(-4163)
$ws.Range("D7").PasteSpecial(-4163)
$ws.Range("E7").PasteSpecial(-4163)
$ws.Range("F7").PasteSpecial(-4163)
$ws.Range("G7").PasteSpecial(-4163)
$ws.Range("I1").Clear()
$excel.CutCopyMode = 0

# Copy id-column (A) style - border + centered alignment - from the last
# existing data row onto the two freshly added id cells.
$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("A7").PasteSpecial(-4122)
$excel.CutCopyMode = 0
